# New grid distribution network calc
# Update the model's grid-generation / distribution-network calculation
# inputs on Sheet1 (row 2) and align AJ1's header formatting with the
# rest of the header row (AI1), matching the new
# "NewGridGenerationCapacityTimestepLimit" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated calculation inputs (row 2) ---
$ws.Range("G2").Value = 1389.601701909612
$ws.Range("H2").Value = 0.1599888484040377
$ws.Range("V2").Value = 0.0986176933596884
$ws.Range("W2").Value = 44.22420440792517
$ws.Range("AA2").Value = 4000
$ws.Range("AB2").Value = 4500

# --- AJ1 header formatting: drop the extra left/right border + fill that
# singled this column out, so it matches the rest of the header row (AI1) ---
$ws.Range("AI1").Copy()
$ws.Range("AJ1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
